$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.221.25"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "2.573.74"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.968.39"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").Value = "2.542.86"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.850"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "43.266.79"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0810"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("E39").Value = "  +4.92%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "1.992.24"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "2.817.07"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "
